$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "NSE:AMBIKCO"
$ws.Range("C2").Value = "NSE:ADANIGREEN"
$ws.Range("D2").Value = "NSE:AXISBANK"
$ws.Range("E2").Value = "NSE:ADANIGREEN"
$ws.Range("F2").Value = "NSE:DLF"

# Row 3
$ws.Range("B3").Value = "NSE:ANANTRAJ"
$ws.Range("C3").Value = "NSE:ASHOKLEY"

# Row 4
$ws.Range("B4").Value = "NSE:APCOTEXIND"
$ws.Range("C4").Value = "NSE:BFINVEST"

# Row 5
$ws.Range("C5").Value = "NSE:BLISSGVS"

# Row 6
$ws.Range("B6").Value = "NSE:ASAHIINDIA"
$ws.Range("C6").Value = "NSE:BLUEJET"

# Row 7
$ws.Range("B7").Value = "NSE:BLBLIMITED"
$ws.Range("C7").Value = "NSE:DHRUV"

# Row 8
$ws.Range("B8").Value = "NSE:BORORENEW"
$ws.Range("C8").Value = "NSE:EPL"

# Row 9
$ws.Range("B9").Value = "NSE:CESC"
$ws.Range("C9").Value = "NSE:GODFRYPHLP"

# Row 10
$ws.Range("B10").Value = "NSE:DLF"
$ws.Range("C10").Value = "NSE:GRINDWELL"

# Row 11
$ws.Range("B11").Value = "NSE:DPSCLTD"
$ws.Range("C11").Value = "NSE:INGERRAND"

# Row 12
$ws.Range("B12").Value = "NSE:ELGIRUBCO"
$ws.Range("C12").Value = "NSE:KANPRPLA"

# Row 13
$ws.Range("B13").Value = "NSE:GMDCLTD"
$ws.Range("C13").Value = "NSE:KIRLOSBROS"

# Row 14
$ws.Range("B14").Value = "NSE:GSS"
$ws.Range("C14").Value = "NSE:NDTV"

# Row 15
$ws.Range("B15").Value = "NSE:GUJAPOLLO"
$ws.Range("C15").Value = "NSE:PDMJEPAPER"

# Row 16 - C16 becomes empty
$ws.Range("B16").Value = "NSE:GUJRAFFIA"
$ws.Range("C16").Value = ""

# Row 17 - C17 becomes empty
$ws.Range("B17").Value = "NSE:JHS"
$ws.Range("C17").Value = ""

# Row 18 - C18 becomes empty
$ws.Range("B18").Value = "NSE:KAMATHOTEL"
$ws.Range("C18").Value = ""

# Row 19 - C19 becomes empty
$ws.Range("B19").Value = "NSE:MCX"
$ws.Range("C19").Value = ""

# Row 20 - C20 becomes empty
$ws.Range("B20").Value = "NSE:MOHITIND"
$ws.Range("C20").Value = ""

# Row 21 - C21 becomes empty
$ws.Range("B21").Value = "NSE:MOMENTUM"
$ws.Range("C21").Value = ""

# Row 22
$ws.Range("B22").Value = "NSE:MTNL"

# Row 23
$ws.Range("B23").Value = "NSE:PARACABLES"

# Row 24
$ws.Range("B24").Value = "NSE:PKTEA"

# Row 25
$ws.Range("B25").Value = "NSE:RAIN"

# Row 26
$ws.Range("B26").Value = "NSE:RELIGARE"

# Row 27
$ws.Range("B27").Value = "NSE:RHL"

# Row 28
$ws.Range("B28").Value = "NSE:RTNPOWER"

# Row 29
$ws.Range("B29").Value = "NSE:SALZERELEC"

# Remove the now-obsolete rows 30-33 (content shifted into rows 2-29 above)
$ws.Rows("30:33").Delete()
